# Auto-generated edit script: updates computed price/profit columns (H-N)
# on several sheets, per scheduled-runner refresh of Golem_Profits data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 866.3333
$ws.Range("I80").Value = 466.1111
$ws.Range("J80").Value = 1266.5555
$ws.Range("K80").Value = 1398.3333
$ws.Range("L80").Value = 3799.6665
$ws.Range("M80").Value = -400.3333
$ws.Range("N80").Value = -5795.666499999999
$ws.Range("H83").Value = 866.3333
$ws.Range("I83").Value = 466.1111
$ws.Range("J83").Value = 1266.5555
$ws.Range("K83").Value = 4194.9999
$ws.Range("L83").Value = 11398.9995
$ws.Range("M83").Value = 797.0001000000002
$ws.Range("N83").Value = -21382.9995
$ws.Range("H92").Value = 83334200
$ws.Range("J92").Value = 1039.8
$ws.Range("L92").Value = 1039.8
$ws.Range("N92").Value = -3535.8
$ws.Range("H132").Value = 2413.7693
$ws.Range("I132").Value = 2805.4
$ws.Range("J132").Value = 1108.3334
$ws.Range("K132").Value = 8416.200000000001
$ws.Range("L132").Value = 3325.0002
$ws.Range("M132").Value = -5886.200000000001
$ws.Range("N132").Value = -8385.0002
$ws.Range("H137").Value = 2250

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2070.3
$ws.Range("I2").Value = 2280.3333
$ws.Range("K2").Value = 2280.3333
$ws.Range("M2").Value = -2167.3333
$ws.Range("H61").Value = 3525
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3424
$ws.Range("H76").Value = 43162.668
$ws.Range("J76").Value = 43162.668
$ws.Range("L76").Value = 43162.668
$ws.Range("N76").Value = -43838.668
$ws.Range("H79").Value = 43162.668
$ws.Range("J79").Value = 43162.668
$ws.Range("L79").Value = 43162.668
$ws.Range("N79").Value = -45502.668
$ws.Range("H97").Value = 30304624
$ws.Range("I97").Value = 30304624
$ws.Range("K97").Value = 30304624
$ws.Range("M97").Value = -30304128
$ws.Range("H116").Value = 2070.3
$ws.Range("I116").Value = 2280.3333
$ws.Range("K116").Value = 2280.3333
$ws.Range("M116").Value = 13.66670000000022
$ws.Range("H122").Value = 7813.5293
$ws.Range("I122").Value = 4305.3335
$ws.Range("J122").Value = 9727.091
$ws.Range("K122").Value = 12916.0005
$ws.Range("L122").Value = 29181.273
$ws.Range("M122").Value = -10466.0005
$ws.Range("N122").Value = -34081.273
$ws.Range("H132").Value = 3007
$ws.Range("I132").Value = 3007
$ws.Range("K132").Value = 9021
$ws.Range("M132").Value = -6491
$ws.Range("H136").Value = 3525
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2070.3
$ws.Range("I3").Value = 2280.3333
$ws.Range("K3").Value = 2280.3333
$ws.Range("M3").Value = -2166.3333
$ws.Range("H22").Value = 449.6154
$ws.Range("I22").Value = 486.91666
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 486.91666
$ws.Range("L22").Value = 2
$ws.Range("M22").Value = -313.91666
$ws.Range("N22").Value = -348
$ws.Range("H36").Value = 13261.6
$ws.Range("I36").Value = 11702.125
$ws.Range("K36").Value = 11702.125
$ws.Range("M36").Value = -11168.125
$ws.Range("H107").Value = 103625
$ws.Range("I107").Value = 134833.33
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 134833.33
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -132913.33
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 561.93335
$ws.Range("I22").Value = 603.0769
$ws.Range("J22").Value = 294.5
$ws.Range("K22").Value = 603.0769
$ws.Range("L22").Value = 294.5
$ws.Range("M22").Value = -253.0769
$ws.Range("N22").Value = -994.5
$ws.Range("H74").Value = 45499.875
$ws.Range("J74").Value = 49499.75
$ws.Range("L74").Value = 49499.75
$ws.Range("N74").Value = -51247.75
$ws.Range("H77").Value = 45499.875
$ws.Range("J77").Value = 49499.75
$ws.Range("L77").Value = 148499.25
$ws.Range("N77").Value = -157235.25
$ws.Range("H99").Value = 4959.8
$ws.Range("I99").Value = 4824.75
$ws.Range("K99").Value = 4824.75
$ws.Range("M99").Value = -3326.75
$ws.Range("H107").Value = 782.1539
$ws.Range("I107").Value = 518.7778
$ws.Range("K107").Value = 518.7778
$ws.Range("M107").Value = 1401.2222
$ws.Range("H126").Value = 4959.8
$ws.Range("I126").Value = 4824.75
$ws.Range("K126").Value = 14474.25
$ws.Range("M126").Value = -12004.25
$ws.Range("H132").Value = 1002.75
$ws.Range("I132").Value = 1002.75
$ws.Range("K132").Value = 3008.25
$ws.Range("M132").Value = -478.25
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35938.207
$ws.Range("I4").Value = 41208.32
$ws.Range("K4").Value = 123624.96
$ws.Range("M4").Value = -123512.96
$ws.Range("H44").Value = 301.46155
$ws.Range("I44").Value = 260.45456
$ws.Range("J44").Value = 527
$ws.Range("K44").Value = 781.36368
$ws.Range("L44").Value = 1581
$ws.Range("M44").Value = -383.36368
$ws.Range("N44").Value = -2377
$ws.Range("H50").Value = 400.5
$ws.Range("I50").Value = 400.5
$ws.Range("K50").Value = 1201.5
$ws.Range("M50").Value = -720.5
$ws.Range("H53").Value = 400.5
$ws.Range("I53").Value = 400.5
$ws.Range("K53").Value = 1201.5
$ws.Range("M53").Value = -720.5
$ws.Range("H80").Value = 2084.1667
$ws.Range("I80").Value = 1849.5
$ws.Range("J80").Value = 2201.5
$ws.Range("K80").Value = 5548.5
$ws.Range("L80").Value = 6604.5
$ws.Range("M80").Value = -4612.5
$ws.Range("N80").Value = -8476.5
$ws.Range("H83").Value = 2084.1667
$ws.Range("I83").Value = 1849.5
$ws.Range("J83").Value = 2201.5
$ws.Range("K83").Value = 16645.5
$ws.Range("L83").Value = 19813.5
$ws.Range("M83").Value = -11965.5
$ws.Range("N83").Value = -29173.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 5249.8
$ws.Range("J22").Value = 5249.8
$ws.Range("L22").Value = 5249.8
$ws.Range("N22").Value = -6307.8
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -29984
$ws.Range("H102").Value = 2318.375
$ws.Range("I102").Value = 999
$ws.Range("J102").Value = 2758.1667
$ws.Range("K102").Value = 999
$ws.Range("L102").Value = 2758.1667
$ws.Range("M102").Value = 623
$ws.Range("N102").Value = -6002.1667
$ws.Range("H107").Value = 83335580
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 111113450
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 111113450
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -111117290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1755
$ws.Range("I10").Value = 956
$ws.Range("K10").Value = 956
$ws.Range("M10").Value = -816
$ws.Range("H74").Value = 82000
$ws.Range("J74").Value = 82000
$ws.Range("L74").Value = 82000
$ws.Range("N74").Value = -83996
$ws.Range("H77").Value = 82000
$ws.Range("J77").Value = 82000
$ws.Range("L77").Value = 246000
$ws.Range("N77").Value = -255984
$ws.Range("H93").Value = 18523692
$ws.Range("I93").Value = 23814196
$ws.Range("J93").Value = 6924.5
$ws.Range("K93").Value = 23814196
$ws.Range("L93").Value = 6924.5
$ws.Range("M93").Value = -23812948
$ws.Range("N93").Value = -9420.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 5618.5713
$ws.Range("I43").Value = 721.6667
$ws.Range("K43").Value = 721.6667
$ws.Range("M43").Value = -572.6667
